$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Files-tab Neo4j query (B4): drop the File Type and Breed
# columns from the RETURN clause.
$ws.Range("B4").Value = "`nMATCH (f:file)-->(parent)`nWITH DISTINCT f, parent`nMATCH (f)-[*]->(c:case)<--(demo:demographic)`n MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`nWHERE diag.stage_of_disease IN ['Vb']`nWITH DISTINCT f, parent, c, demo, diag, s`nRETURN coalesce(f.file_name, '') AS ``File Name``, `n        coalesce(labels(parent)[0], '') AS ``Association``,`n        coalesce(f.file_description, '') AS ``Description``,`n        coalesce(f.file_format, '') AS ``Format``,`n        coalesce(f.file_size, '') AS ``Size``,`n        coalesce(c.case_id, '') AS ``Case ID``, `n        coalesce(diag.disease_term,'') AS Diagnosis , `n        coalesce(s.clinical_study_designation,'') AS ``Study Code``"

# Move the selection/scroll position to match the author's saved view
# (scrolled down to row 4, cell B4 selected).
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("B4").Select()
